$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# E2: 42.3 -> 46
$ws.Range("E2").Value = 46

# Row 10 updates: C10 758 -> 227.4, D10 formula(=D3) -> plain 42, E10 formula(=E3) -> plain 44
$ws.Range("C10").Value = 227.4
$ws.Range("D10").Value = 42
$ws.Range("E10").Value = 44

# Update the sheet's active selection (and drop the scrolled topLeftCell) to D13
$ws.Range("D13").Select()
